$wb = $excel.ActiveWorkbook

$wsList = $wb.Worksheets.Item("LIST")
$wsFeuil1 = $wb.Worksheets.Item("Feuil1")

# Update the TODO LIST ("LIST" sheet): shrink the list of test cases down to
# the first one, and add the new test-result entry ("RESULTAT DES TESTS").
$wsList.Range("A2").Value = "AD.SEC.001.FON.01"
$wsList.Range("A3").Value = "RO.FOU.001.SUP.01"
$wsList.Range("A3").ClearFormats()
$wsList.Range("A4").Value = ""
$wsList.Range("A5").Value = ""
$wsList.Range("A6").Value = ""
$wsList.Range("A7").Value = ""
$wsList.Range("A8").Value = ""

# Restore the saved view state (selections) on each sheet, making sure
# "LIST" ends up as the active/selected tab again.
$wsFeuil1.Range("F2").Select()
$wsList.Activate()
$wsList.Range("C10").Select()
